# Update "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps, as produced by a fresh
# "Generate Report for Handback" run.
#
# Overview!G2               2016-09-04 07:11:10 -> 2016-09-04 07:12:29
# zh-cn!H2  (Handoff)       2016-09-04 07:11:01 -> 2016-09-04 07:12:25
# zh-cn!K2  (Handback)      2016-09-04 07:11:59 -> 2016-09-04 07:12:43
# de-de!H2  (Handoff)       2016-09-04 07:11:10 -> 2016-09-04 07:12:29
# de-de!K2  (Handback)      2016-09-04 07:12:12 -> 2016-09-04 07:12:50

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$overview.Range("G2").Value = "2016-09-04 07:12:29"

$zhcn.Range("H2").Value = "2016-09-04 07:12:25"
$zhcn.Range("K2").Value = "2016-09-04 07:12:43"

$dede.Range("H2").Value = "2016-09-04 07:12:29"
$dede.Range("K2").Value = "2016-09-04 07:12:50"
